# Auto-generated Excel COM-interop script
# Applies literal value changes per the target diff for sheet1 (Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated / added cell values
$ws.Cells.Item(2, 2).Value = 2.411890494488538
$ws.Cells.Item(2, 3).Value = 10.32257956503014
$ws.Cells.Item(2, 4).Value = -7.980392560045654
$ws.Cells.Item(2, 5).Value = 0.1714771788446484
$ws.Cells.Item(2, 6).Value = 1.965700931935187
$ws.Cells.Item(2, 7).Value = -1.209548186156553
$ws.Cells.Item(2, 8).Value = -1.50346969825023
$ws.Cells.Item(2, 9).Value = 1.011663206852835
$ws.Cells.Item(2, 10).Value = 0.08227779023355938
$ws.Cells.Item(2, 11).Value = 0.3174504178062259
$ws.Cells.Item(3, 2).Value = 9.407028267080065
$ws.Cells.Item(3, 3).Value = -8.535831474695623
$ws.Cells.Item(3, 4).Value = -0.1474425313531172
$ws.Cells.Item(3, 5).Value = 1.756498547467592
$ws.Cells.Item(3, 6).Value = -1.36276561796536
$ws.Cells.Item(3, 7).Value = -1.628924567856255
$ws.Cells.Item(3, 8).Value = 0.9000908150593795
$ws.Cells.Item(3, 9).Value = -0.02236996370037081
$ws.Cells.Item(3, 10).Value = 0.2162591392419525
$ws.Cells.Item(3, 11).Value = -0.08907140834264399
$ws.Cells.Item(4, 2).Value = -15.70106135728257
$ws.Cells.Item(4, 3).Value = -5.922097328900056
$ws.Cells.Item(4, 4).Value = -2.885559413464426
$ws.Cells.Item(4, 5).Value = -5.113293446597504
$ws.Cells.Item(4, 6).Value = -4.685009736670638
$ws.Cells.Item(4, 7).Value = -1.619613956105766
$ws.Cells.Item(4, 8).Value = -2.12984807464721
$ws.Cells.Item(4, 9).Value = -1.575446602675451
$ws.Cells.Item(4, 10).Value = -1.639399916244841
$ws.Cells.Item(4, 11).Value = -0.9435672860403349
$ws.Cells.Item(5, 2).Value = -4.241916182330832
$ws.Cells.Item(5, 3).Value = 4.530485349044278
$ws.Cells.Item(5, 4).Value = -2.865237130372411
$ws.Cells.Item(5, 5).Value = -0.703648721481464
$ws.Cells.Item(5, 6).Value = 0.3585679458042089
$ws.Cells.Item(5, 7).Value = 0.2821581714231107
$ws.Cells.Item(5, 8).Value = 0.01456377582551255
$ws.Cells.Item(5, 9).Value = 0.003172630064957882
$ws.Cells.Item(5, 10).Value = 0.3397553536669872
$ws.Cells.Item(5, 11).Value = 0.1548536679346297
$ws.Cells.Item(6, 2).Value = 0.8032788223262832
$ws.Cells.Item(6, 3).Value = -1.002366475608788
$ws.Cells.Item(6, 4).Value = -1.464557645652444
$ws.Cells.Item(6, 5).Value = 0.9087989414264609
$ws.Cells.Item(6, 6).Value = 0.03490493049987109
$ws.Cells.Item(6, 7).Value = 0.2751479164465338
$ws.Cells.Item(6, 8).Value = -0.04060839016473031
$ws.Cells.Item(6, 9).Value = 0.4714428511361884
$ws.Cells.Item(6, 10).Value = 0.1848835897783218
$ws.Cells.Item(6, 11).Value = 0.2125075656625323
$ws.Cells.Item(7, 2).Value = -0.5941340219870179
$ws.Cells.Item(7, 3).Value = -1.21293460176605
$ws.Cells.Item(7, 4).Value = 0.8268928876033257
$ws.Cells.Item(7, 5).Value = 0.1099969798811567
$ws.Cells.Item(7, 6).Value = 0.3810550769602576
$ws.Cells.Item(7, 7).Value = 0.0178082559939251
$ws.Cells.Item(7, 8).Value = 0.5378815349662799
$ws.Cells.Item(7, 9).Value = 0.2595655483891583
$ws.Cells.Item(7, 10).Value = 0.2828208575635111
$ws.Cells.Item(7, 11).Value = 0.3749895042266514
$ws.Cells.Item(8, 2).Value = -1.318433813614865
$ws.Cells.Item(8, 3).Value = 0.7007876036379678
$ws.Cells.Item(8, 4).Value = 0.2782450485143884
$ws.Cells.Item(8, 5).Value = 0.4417191695641399
$ws.Cells.Item(8, 6).Value = 0.02894109536855799
$ws.Cells.Item(8, 7).Value = 0.6007442862932105
$ws.Cells.Item(8, 8).Value = 0.3169170761829015
$ws.Cells.Item(8, 9).Value = 0.3273260345678901
$ws.Cells.Item(8, 10).Value = 0.42606337525307
$ws.Cells.Item(8, 11).Value = 0.6412619431822899
$ws.Cells.Item(9, 2).Value = 0.0963469837902291
$ws.Cells.Item(9, 3).Value = 0.214041671159695
$ws.Cells.Item(9, 4).Value = 0.7332664776213567
$ws.Cells.Item(9, 5).Value = -0.008765780717375604
$ws.Cells.Item(9, 6).Value = 0.5921602662197494
$ws.Cells.Item(9, 7).Value = 0.4006581537802698
$ws.Cells.Item(9, 8).Value = 0.3621531794959351
$ws.Cells.Item(9, 9).Value = 0.452627758305367
$ws.Cells.Item(9, 10).Value = 0.6867607839288887
$ws.Cells.Item(9, 11).Value = 0.1781578843816368
$ws.Cells.Item(10, 2).Value = 0.1264008423207837
$ws.Cells.Item(10, 3).Value = 0.6808472755916881
$ws.Cells.Item(10, 4).Value = 0.04988061626763002
$ws.Cells.Item(10, 5).Value = 0.5970894115568507
$ws.Cells.Item(10, 6).Value = 0.3895648707313746
$ws.Cells.Item(10, 7).Value = 0.3749198787210216
$ws.Cells.Item(10, 8).Value = 0.4613573173527261
$ws.Cells.Item(10, 9).Value = 0.6892627280777406
$ws.Cells.Item(10, 10).Value = 0.1844338218533179
$ws.Cells.Item(10, 11).Value = 0.4715052544735016
$ws.Cells.Item(11, 2).Value = 0.7032752552246967
$ws.Cells.Item(11, 3).Value = 0.0524841558300787
$ws.Cells.Item(11, 4).Value = 0.5750217259028355
$ws.Cells.Item(11, 5).Value = 0.3835828904270196
$ws.Cells.Item(11, 6).Value = 0.3695766285386105
$ws.Cells.Item(11, 7).Value = 0.4509887839823598
$ws.Cells.Item(11, 8).Value = 0.6807854059541167
$ws.Cells.Item(11, 9).Value = 0.1766615135465071
$ws.Cells.Item(11, 10).Value = 0.4628812809405329
$ws.Cells.Item(11, 11).Value = 0.2074249537672726
$ws.Cells.Item(12, 2).Value = 0.0739087272872988
$ws.Cells.Item(12, 3).Value = 0.6951995747020479
$ws.Cells.Item(12, 4).Value = 0.2943885460132365
$ws.Cells.Item(12, 5).Value = 0.3398496246900327
$ws.Cells.Item(12, 6).Value = 0.4596149040122699
$ws.Cells.Item(12, 7).Value = 0.6548452325286815
$ws.Cells.Item(12, 8).Value = 0.1538829505182796
$ws.Cells.Item(12, 9).Value = 0.4492632457919151
$ws.Cells.Item(12, 10).Value = 0.1889041565820968
$ws.Cells.Item(12, 11).Value = 0.5706702220727796
$ws.Cells.Item(13, 2).Value = 0.6528789423816584
$ws.Cells.Item(13, 3).Value = 0.2651840721575033
$ws.Cells.Item(13, 4).Value = 0.3368290248851115
$ws.Cells.Item(13, 5).Value = 0.4407536204007895
$ws.Cells.Item(13, 6).Value = 0.6345141014634773
$ws.Cells.Item(13, 7).Value = 0.1393527950840318
$ws.Cells.Item(13, 8).Value = 0.4326222002996472
$ws.Cells.Item(13, 9).Value = 0.1713203111533466
$ws.Cells.Item(13, 10).Value = 0.5541963385427369
$ws.Cells.Item(13, 11).Value = 0.1437698493309027
$ws.Cells.Item(14, 2).Value = 0.6053818127754134
$ws.Cells.Item(14, 3).Value = 0.4122000866690486
$ws.Cells.Item(14, 4).Value = 0.2496603340877904
$ws.Cells.Item(14, 5).Value = 0.6592080140502106
$ws.Cells.Item(14, 6).Value = 0.1487321986403278
$ws.Cells.Item(14, 7).Value = 0.3778114016882561
$ws.Cells.Item(14, 8).Value = 0.1524262202646768
$ws.Cells.Item(14, 9).Value = 0.5393323377276911
$ws.Cells.Item(14, 10).Value = 0.115058138701532
$ws.Cells.Item(14, 11).Value = 0.4067718394308724
$ws.Cells.Item(15, 2).Value = 0.8628949586592991
$ws.Cells.Item(15, 3).Value = 0.2967710363001488
$ws.Cells.Item(15, 4).Value = 0.4189247832594023
$ws.Cells.Item(15, 5).Value = 0.1846772797061906
$ws.Cells.Item(15, 6).Value = 0.3757606442486632
$ws.Cells.Item(15, 7).Value = 0.07585798082864662
$ws.Cells.Item(15, 8).Value = 0.5119329433524077
$ws.Cells.Item(15, 9).Value = 0.08876908850380663
$ws.Cells.Item(15, 10).Value = 0.3633745487175398
$ws.Cells.Item(16, 2).Value = 0.6090966232236873
$ws.Cells.Item(16, 3).Value = 0.5522135229949265
$ws.Cells.Item(16, 4).Value = 0.005598857889999004
$ws.Cells.Item(16, 5).Value = 0.4039548830192304
$ws.Cells.Item(16, 6).Value = 0.1118832920210401
$ws.Cells.Item(16, 7).Value = 0.4743913731481941
$ws.Cells.Item(16, 8).Value = 0.08322674941644539
$ws.Cells.Item(16, 9).Value = 0.3675498776562884
$ws.Cells.Item(17, 2).Value = 0.7878040141027678
$ws.Cells.Item(17, 3).Value = 0.09027759876430858
$ws.Cells.Item(17, 4).Value = 0.2583545163855133
$ws.Cells.Item(17, 5).Value = 0.1218370348802827
$ws.Cells.Item(17, 6).Value = 0.489756542847739
$ws.Cells.Item(17, 7).Value = 0.04437841445902233
$ws.Cells.Item(17, 8).Value = 0.3517040686291025
$ws.Cells.Item(18, 2).Value = 0.4013017852456914
$ws.Cells.Item(18, 3).Value = 0.3754432907967085
$ws.Cells.Item(18, 4).Value = -0.04062710656928412
$ws.Cells.Item(18, 5).Value = 0.5187154933129405
$ws.Cells.Item(18, 6).Value = 0.08012128691392592
$ws.Cells.Item(18, 7).Value = 0.3203764222454754
$ws.Cells.Item(19, 2).Value = 0.6222684682008229
$ws.Cells.Item(19, 3).Value = -0.02297123903139461
$ws.Cells.Item(19, 4).Value = 0.4240932542019461
$ws.Cells.Item(19, 5).Value = 0.112338675162406
$ws.Cells.Item(19, 6).Value = 0.3327645480731927
$ws.Cells.Item(20, 2).Value = 0.2167051203848173
$ws.Cells.Item(20, 3).Value = 0.5091174976711597
$ws.Cells.Item(20, 4).Value = -0.004145903195608092
$ws.Cells.Item(20, 5).Value = 0.3478698197250452
$ws.Cells.Item(21, 2).Value = 0.6739775747052469
$ws.Cells.Item(21, 3).Value = 0.009391369052308113
$ws.Cells.Item(21, 4).Value = 0.2848969007350822
$ws.Cells.Item(22, 2).Value = 0.2632404109177161
$ws.Cells.Item(22, 3).Value = 0.3842149509171186
$ws.Cells.Item(23, 2).Value = 0.4282746421565676

# Cells removed in the target (clear contents)
$ws.Cells.Item(15, 11).ClearContents()
$ws.Cells.Item(16, 10).ClearContents()
$ws.Cells.Item(17, 9).ClearContents()
$ws.Cells.Item(18, 8).ClearContents()
$ws.Cells.Item(19, 7).ClearContents()
$ws.Cells.Item(20, 6).ClearContents()
$ws.Cells.Item(21, 5).ClearContents()
$ws.Cells.Item(22, 4).ClearContents()
$ws.Cells.Item(23, 3).ClearContents()
$ws.Cells.Item(24, 2).ClearContents()
